# The wedding-planner "Dataset" sheet is being tidied up in real Excel:
# the FamilyGroup column (column B) is widened so the longer group
# names (e.g. "Anderson Family", "Work Colleagues") are fully visible.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).ColumnWidth = 21.1
